$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data values (order matters for shared-string slot reuse:
# set C2 first, then B2, then E2, so the rebuilt sharedStrings table keeps
# the same index layout as before the edit)
$ws.Range("A2").Value = 204
$ws.Range("C2").Value = "mdex:double"
$ws.Range("B2").Value = "BASE_TRANSACTION_VALUE"
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = "Base Transaction Value"

# Set new column widths for columns F and G (closest values the host's
# pixel-quantized ColumnWidth model can represent; target stored widths
# are 17.140625 / 12.7109375 but the engine snaps to 1/6-character steps)
$ws.Columns.Item(6).ColumnWidth = 17.04
$ws.Columns.Item(7).ColumnWidth = 11.85

# Update selection to D1
[void]$ws.Range("D1").Select()
